$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh cryptocurrency price/volume snapshot (GitHub Actions data pull) ---

# Cells whose new text looks like a plain number (e.g. "1.005") must be forced
# to Text first, otherwise Excel auto-converts the input into a numeric value.
$textForcedCells = @(
    @(4, 4),
    @(5, 4),
    @(7, 4),
    @(8, 4),
    @(9, 4),
    @(10, 4),
    @(11, 4),
    @(12, 4),
    @(13, 4),
    @(14, 4),
    @(15, 4),
    @(16, 4),
    @(18, 4),
    @(19, 4),
    @(20, 4),
    @(21, 4),
    @(22, 4),
    @(23, 4),
    @(25, 4),
    @(26, 4),
    @(27, 4),
    @(28, 4),
    @(29, 4),
    @(30, 4),
    @(31, 4),
    @(32, 4),
    @(34, 4),
    @(35, 4),
    @(36, 4),
    @(37, 4),
    @(38, 4),
    @(39, 4),
    @(40, 4),
    @(41, 4),
    @(42, 4),
    @(43, 4),
    @(44, 4),
    @(45, 4),
    @(46, 4),
    @(47, 4),
    @(49, 4),
    @(50, 4),
    @(51, 4)
)

foreach ($coord in $textForcedCells) {
    $ws.Cells.Item($coord[0], $coord[1]).NumberFormat = "@"
}

# Apply the updated values
$ws.Cells.Item(4, 4).Value = '1.005'
$ws.Cells.Item(5, 4).Value = '308.68'
$ws.Cells.Item(7, 4).Value = '0.3916'
$ws.Cells.Item(8, 4).Value = '0.3869'
$ws.Cells.Item(9, 4).Value = '1.001'
$ws.Cells.Item(10, 4).Value = '1.367'
$ws.Cells.Item(11, 4).Value = '49.26'
$ws.Cells.Item(12, 4).Value = '0.08584'
$ws.Cells.Item(13, 4).Value = '23.55'
$ws.Cells.Item(14, 4).Value = '7.089'
$ws.Cells.Item(15, 4).Value = '0.00001283'
$ws.Cells.Item(16, 4).Value = '7.509'
$ws.Cells.Item(18, 4).Value = '94.16'
$ws.Cells.Item(19, 4).Value = '0.06912'
$ws.Cells.Item(20, 4).Value = '20.29'
$ws.Cells.Item(21, 4).Value = '6.901'
$ws.Cells.Item(22, 4).Value = '1.002'
$ws.Cells.Item(23, 4).Value = '13.58'
$ws.Cells.Item(25, 4).Value = '2.402'
$ws.Cells.Item(26, 4).Value = '2.861'
$ws.Cells.Item(27, 4).Value = '22.32'
$ws.Cells.Item(28, 4).Value = '158.13'
$ws.Cells.Item(29, 4).Value = '140.59'
$ws.Cells.Item(30, 4).Value = '8.175'
$ws.Cells.Item(31, 4).Value = '5.272'
$ws.Cells.Item(32, 4).Value = '2.477'
$ws.Cells.Item(34, 4).Value = '0.08126'
$ws.Cells.Item(35, 4).Value = '6.805'
$ws.Cells.Item(36, 4).Value = '0.02904'
$ws.Cells.Item(37, 4).Value = '0.9612'
$ws.Cells.Item(38, 4).Value = '0.2686'
$ws.Cells.Item(39, 4).Value = '0.09187'
$ws.Cells.Item(40, 4).Value = '10.27'
$ws.Cells.Item(41, 4).Value = '1.440'
$ws.Cells.Item(42, 4).Value = '0.7510'
$ws.Cells.Item(43, 4).Value = '13.02'
$ws.Cells.Item(44, 4).Value = '16.13'
$ws.Cells.Item(45, 4).Value = '0.6890'
$ws.Cells.Item(46, 4).Value = '2.461'
$ws.Cells.Item(47, 4).Value = '4.098'
$ws.Cells.Item(49, 4).Value = '0.08368'
$ws.Cells.Item(50, 4).Value = '1.265'
$ws.Cells.Item(51, 4).Value = '133.56'
$ws.Cells.Item(2, 4).Value = '24.062.15'
$ws.Cells.Item(2, 5).Value = '  -2.64%  '
$ws.Cells.Item(3, 4).Value = '1.646.04'
$ws.Cells.Item(3, 5).Value = '  -2.05%  '
$ws.Cells.Item(4, 5).Value = '  +0.28%  '
$ws.Cells.Item(5, 5).Value = '  -1.56%  '
$ws.Cells.Item(6, 5).Value = '  -0.02%  '
$ws.Cells.Item(7, 5).Value = '  -0.43%  '
$ws.Cells.Item(8, 5).Value = '  -2.40%  '
$ws.Cells.Item(9, 5).Value = '  -0.05%  '
$ws.Cells.Item(10, 5).Value = '  -3.54%  '
$ws.Cells.Item(11, 5).Value = '  -5.01%  '
$ws.Cells.Item(12, 5).Value = '  -1.02%  '
$ws.Cells.Item(13, 5).Value = '  -7.21%  '
$ws.Cells.Item(14, 5).Value = '  -3.66%  '
$ws.Cells.Item(15, 5).Value = '  -3.08%  '
$ws.Cells.Item(16, 5).Value = '  -3.94%  '
$ws.Cells.Item(17, 4).Value = '1.648.84'
$ws.Cells.Item(17, 5).Value = '  -2.22%  '
$ws.Cells.Item(18, 5).Value = '  +0.31%  '
$ws.Cells.Item(19, 5).Value = '  -2.76%  '
$ws.Cells.Item(20, 5).Value = '  +0.02%  '
$ws.Cells.Item(21, 5).Value = '  -3.31%  '
$ws.Cells.Item(22, 5).Value = '  -0.20%  '
$ws.Cells.Item(23, 5).Value = '  -3.38%  '
$ws.Cells.Item(24, 4).Value = '24.074.84'
$ws.Cells.Item(24, 5).Value = '  -2.61%  '
$ws.Cells.Item(25, 5).Value = '  +2.10%  '
$ws.Cells.Item(26, 5).Value = '  +2.77%  '
$ws.Cells.Item(27, 5).Value = '  -6.13%  '
$ws.Cells.Item(28, 5).Value = '  -3.01%  '
$ws.Cells.Item(29, 2).Value = 'BitcoinCash'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(29, 5).Value = '  -7.23%  '
$ws.Cells.Item(30, 2).Value = 'Filecoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(30, 5).Value = '  +3.59%  '
$ws.Cells.Item(31, 2).Value = 'HuobiToken'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(31, 5).Value = '  -9.15%  '
$ws.Cells.Item(32, 5).Value = '  +4.49%  '
$ws.Cells.Item(33, 4).Value = '1.829.16'
$ws.Cells.Item(33, 5).Value = '  -2.15%  '
$ws.Cells.Item(34, 5).Value = '  -4.22%  '
$ws.Cells.Item(35, 5).Value = '  -2.59%  '
$ws.Cells.Item(36, 5).Value = '  -6.41%  '
$ws.Cells.Item(37, 5).Value = '  -5.40%  '
$ws.Cells.Item(38, 5).Value = '  -4.30%  '
$ws.Cells.Item(39, 5).Value = '  -3.41%  '
$ws.Cells.Item(40, 5).Value = '  -2.73%  '
$ws.Cells.Item(41, 5).Value = '  -3.07%  '
$ws.Cells.Item(42, 5).Value = '  -5.86%  '
$ws.Cells.Item(43, 5).Value = '  -4.82%  '
$ws.Cells.Item(44, 5).Value = '  -3.57%  '
$ws.Cells.Item(45, 5).Value = '  -3.87%  '
$ws.Cells.Item(46, 5).Value = '  -4.79%  '
$ws.Cells.Item(47, 5).Value = '  -1.80%  '
$ws.Cells.Item(48, 5).Value = '  -0.05%  '
$ws.Cells.Item(49, 5).Value = '  -3.72%  '
$ws.Cells.Item(50, 5).Value = '  -5.36%  '
$ws.Cells.Item(51, 5).Value = '  -3.70%  '

# Restore default (General/Normal) cell style now that the text value is locked in,
# so no stray number-format styling is left behind on these cells.
foreach ($coord in $textForcedCells) {
    $ws.Cells.Item($coord[0], $coord[1]).Style = "Normal"
}
